$d = $word.ActiveDocument

# Both tables in the document have 6 columns; columns 4 and 5 need their
# widths nudged by 1 twentieth of a point (dxa): column 4 grows from its
# current width to 1515 dxa, and column 5 shrinks to 1508 dxa.
# Setting a single cell's Width updates the whole column (tcW in every row
# plus the corresponding gridCol entry), which is what we want here.

for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $t = $d.Tables.Item($i)

    $cell4 = $t.Cell(1, 4)
    $cell4.Width = 1515 / 20

    $cell5 = $t.Cell(1, 5)
    $cell5.Width = 1508 / 20
}
